$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ETL: limpieza de datos ---
# Delete rows with missing/invalid records (keep ids 1,3,4,7,9,10), from bottom to top
# so row indices above are not shifted while deleting.
$rowsToDelete = @(16,15,14,13,12,9,7,6,3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Fix / normalize remaining dirty values
$ws.Range("D4").Value = "SIN CIUDAD"      # id 4, ciudad Cusco -> SIN CIUDAD
$ws.Range("B6").Value = "Luis Fernandez"  # id 9, remove accent
$ws.Range("B7").Value = "SIN NOMBRE"      # id 10, nombre Patricia Rios -> SIN NOMBRE
